$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date (column C) for rows 2-27 from 45243 to 45244
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = 45244
}
